$d = $word.ActiveDocument

# --- Remove the _GoBack bookmark from its current location (hover paragraph) ---
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# --- Templates used to paste text while controlling the w:rFonts hint="eastAsia" ---
# No-hint template: "ll bar" (6 chars) - run without rFonts hint
$noHintSrc = $d.Range(121, 127)
$ftNoHint = $noHintSrc.FormattedText
$noHintLen = 6

# Hint template: "无效" (2 chars) - run with rFonts hint="eastAsia"
$hintSrc = $d.Range(158, 160)
$ftHint = $hintSrc.FormattedText
$hintLen = 2

function PasteSeg($doc, $pos, $text, $ft, $templateLen) {
    $ph = ""
    for ($i = 0; $i -lt $templateLen; $i++) { $ph = $ph + "X" }
    $r0 = $doc.Range($pos, $pos)
    $r0.InsertAfter($ph)
    $r1 = $doc.Range($pos, $pos + $templateLen)
    $r1.FormattedText = $ft
    $r2 = $doc.Range($pos, $pos + $templateLen)
    $r2.Text = $text
    return $pos + $text.Length
}

# --- New paragraph 1: 为何scroll区域要外套两个应用Internal css的div ---
$lastP = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$lastP.Collapse(0)
$lastP.InsertParagraphAfter()
$newPara1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$cur = $newPara1.Range
$cur.Collapse(1)
$pos = $cur.Start

$pos = PasteSeg $d $pos "为何" $ftHint $hintLen
$pos = PasteSeg $d $pos "scroll" $ftNoHint $noHintLen
$pos = PasteSeg $d $pos "区域" $ftHint $hintLen
$pos = PasteSeg $d $pos "要外套两个应用Internal css" $ftNoHint $noHintLen
$pos = PasteSeg $d $pos "的div" $ftHint $hintLen

# --- New paragraph 2: 主播榜单板块：周榜月榜按钮需要按两次才能触发效果。 ---
$lastP2 = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$lastP2.Collapse(0)
$lastP2.InsertParagraphAfter()
$newPara2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$cur2 = $newPara2.Range
$cur2.Collapse(1)
$pos2 = $cur2.Start

$pos2 = PasteSeg $d $pos2 "主播" $ftHint $hintLen
$pos2 = PasteSeg $d $pos2 "榜单板块：" $ftNoHint $noHintLen
$pos2 = PasteSeg $d $pos2 "周榜月" $ftHint $hintLen
$pos2 = PasteSeg $d $pos2 "榜按钮需要按两次才能触发效果。" $ftNoHint $noHintLen

# --- Re-add the _GoBack bookmark at the very end of the new last paragraph ---
# NOTE: placing a bookmark at the absolute last insertion point of the whole
# document (Content.End - 1) is mishandled by this engine (it silently drops
# the bookmark at document position 0 instead). Work around it by appending a
# throwaway placeholder character, anchoring the bookmark just before it, and
# then deleting the placeholder again - the bookmark stays put.
$placeholderRange = $d.Range($pos2, $pos2)
$placeholderRange.InsertAfter("Z")
$endRange = $d.Range($pos2, $pos2)
$d.Bookmarks.Add("_GoBack", $endRange)
$delRange = $d.Range($pos2, $pos2 + 1)
$delRange.Delete()

Write-Host "Paragraph count:" $d.Paragraphs.Count
Write-Host "New para1:" $newPara1.Range.Text
Write-Host "New para2:" $newPara2.Range.Text
